$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Accounts whose rows must be removed entirely.
$deleteAccounts = @(
    "004452912",  # BRUNO      70205.62
    "004404342",  # ADSON      34005.84
    "000806386",  # FERNANDA   12979.5
    "005631648",  # GLEYSERSON 10000
    "004756968",  # DANIELY    4582.71
    "004231371",  # ADRIANO    500
    "005110867"   # DIG        493.5
)

foreach ($acct in $deleteAccounts) {
    $found = $ws.Columns(1).Find($acct, $null, $null, 1)
    if ($found) {
        $ws.Rows($found.Row).Delete()
    }
}

# Insert a new row for VINICIUS right before the row that holds account 008115265 (ELAINE).
$anchor = $ws.Columns(1).Find("008115265", $null, $null, 1)
$anchorRow = $anchor.Row
$ws.Rows($anchorRow).Insert()

$ws.Cells.Item($anchorRow, 1).Value = "'005886225"
$ws.Cells.Item($anchorRow, 2).Value = "VINICIUS"
$ws.Cells.Item($anchorRow, 3).Value = 5000
